$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B6 currently holds numeric 12 (quantity). Re-enter it as text "12"
# (quote-prefixed, like the existing pincode cell B1) and add a new
# text value "7" in C6 next to it.
$ws.Range("B6").Value = "'12"
$ws.Range("C6").Value = "'7"

$ws.Range("C6").Select()
